$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Header row - author name and week number
$ws.Range("C1").Value = "Jesse Hare"

# Row 3 & 4 - task entries ("Stage" filled first for both rows, then "Task")
$ws.Range("A3").Value = "Project Planning"
$ws.Range("A4").Value = "Project Planning"
$ws.Range("B4").Value = "Write product specification for searcher program"
$ws.Range("B3").Value = "Prepare Project Plan for client, work on introduction"
$ws.Range("C3").Value = "``8"

$ws.Range("E1").Value = 2
$ws.Range("D3").Value = 10
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 10

# Update selection to match final cursor position
$ws.Range("D7").Select()
